# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns
# for the current snapshot of coin data on the active worksheet.
# Price cells that look like plain decimals are forced to Text format
# before assignment so they keep matching the source feed's string
# formatting (e.g. "357.16") instead of being auto-typed as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.069.18'
$ws.Range("E2").Value = '  -0.68%  '
$ws.Range("D3").Value = '2.916.66'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '357.16'
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.76'
$ws.Range("E6").Value = '  -1.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.568'
$ws.Range("E7").Value = '  +1.55%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.635'
$ws.Range("E9").Value = '  +0.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.95'
$ws.Range("E10").Value = '  -2.47%  '
$ws.Range("E11").Value = '  +1.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0869'
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.52'
$ws.Range("E13").Value = '  -1.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.78'
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").Value = '3.375.70'
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("D16").Value = '2.925.60'
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.989'
$ws.Range("E17").Value = '  -1.87%  '
$ws.Range("D18").Value = '52.050.70'
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("E19").Value = '  +4.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.54'
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.89'
$ws.Range("E21").Value = '  -1.69%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.58'
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.21'
$ws.Range("E24").Value = '  -0.40%  '
$ws.Range("E25").Value = '  +1.45%  '
$ws.Range("E26").Value = '  +9.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.67'
$ws.Range("E27").Value = '  +16.05%  '
$ws.Range("E28").Value = '  +0.56%  '
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  +7.96%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.50'
$ws.Range("E31").Value = '  -1.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '37.61'
$ws.Range("E32").Value = '  -0.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.18'
$ws.Range("E33").Value = '  -3.21%  '
$ws.Range("E34").Value = '  -1.41%  '
$ws.Range("E35").Value = '  -2.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0441'
$ws.Range("E36").Value = '  -1.93%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.20'
$ws.Range("E38").Value = '  -3.29%  '
$ws.Range("E39").Value = '  -2.86%  '
$ws.Range("E40").Value = '  -3.47%  '
$ws.Range("E41").Value = '  -4.66%  '
$ws.Range("E42").Value = '  +2.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.90'
$ws.Range("E43").Value = '  -2.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.64'
$ws.Range("E45").Value = '  -0.99%  '
$ws.Range("E46").Value = '  -2.49%  '
$ws.Range("E47").Value = '  -4.63%  '
$ws.Range("D48").Value = '2.129.63'
$ws.Range("E48").Value = '  -3.26%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.254'
$ws.Range("E49").Value = '  -4.84%  '
$ws.Range("E50").Value = '  +0.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.922'
$ws.Range("E51").Value = '  -5.33%  '
